# Updates cryptos list data (price & volume columns), plus a couple of
# row re-orderings (Toncoin/Cardano and SuiNetwork/Dai swapped ranking positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "22.99", "4.78") are not silently coerced into floating
# point numbers when assigned via .Value.
$ws.Range("D2:E51").NumberFormat = "@"

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B -ne $null) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($C -ne $null) { $ws.Cells.Item($Row, 3).Value = $C }
    if ($D -ne $null) { $ws.Cells.Item($Row, 4).Value = $D }
    if ($E -ne $null) { $ws.Cells.Item($Row, 5).Value = $E }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "67.853.59" "  -0.74%  "

# Row 3 - Ethereum
Set-Row 3 $null $null "2.423.94" "  -0.95%  "

# Row 4 - TetherUSD
Set-Row 4 $null $null $null "  -0.03%  "

# Row 5 - BNB
Set-Row 5 $null $null "551.69" "  -0.58%  "

# Row 6 - Solana
Set-Row 6 $null $null "160.56" "  -0.48%  "

# Row 7 - USDC
Set-Row 7 $null $null $null "  -0.01%  "

# Row 8 - XRP
Set-Row 8 $null $null "0.511" "  +2.04%  "

# Row 9 - Dogecoin
Set-Row 9 $null $null $null "  +7.51%  "

# Row 10 - TRON
Set-Row 10 $null $null $null "  -0.33%  "

# Row 11 - was Toncoin, now Cardano
Set-Row 11 "Cardano" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada" "0.326" "  -2.06%  "

# Row 12 - was Cardano, now Toncoin
Set-Row 12 "Toncoin" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton" "4.78" "  +0.16%  "

# Row 13 - WrappedBTC
Set-Row 13 $null $null "67.790.53" "  -0.66%  "

# Row 14 - ShibaInu
Set-Row 14 $null $null "0.0000169" "  +1.20%  "

# Row 15 - Avalanche
Set-Row 15 $null $null "22.99" "  -0.94%  "

# Row 16 - Chainlink
Set-Row 16 $null $null "10.30" "  -3.89%  "

# Row 17 - BitcoinCash
Set-Row 17 $null $null "335.40" "  -1.09%  "

# Row 18 - Uniswap
Set-Row 18 $null $null "6.82" "  -2.11%  "

# Row 19 - Polkadot
Set-Row 19 $null $null $null "  +0.22%  "

# Row 20 - was SuiNetwork, now Dai
Set-Row 20 "Dai" "https://coinranking.com/coin/MoTuySvg7+dai-dai" "1.00" "  +0.02%  "

# Row 21 - was Dai, now SuiNetwork
Set-Row 21 "SuiNetwork" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui" "1.88" "  +2.08%  "

# Row 22 - Litecoin
Set-Row 22 $null $null "66.50" "  +0.36%  "

# Row 23 - NEARProtocol
Set-Row 23 $null $null "3.62" "  -0.66%  "

# Row 24 - Aptos
Set-Row 24 $null $null "8.08" "  +0.68%  "

# Row 25 - PEPE
Set-Row 25 $null $null $null "  +0.19%  "

# Row 26 - InternetComputer(DFINITY)
Set-Row 26 $null $null $null "  +0.05%  "

# Row 27 - FirstDigitalUSD
Set-Row 27 $null $null $null "  +0.02%  "

# Row 28 - Bittensor
Set-Row 28 $null $null "420.84" "  -2.97%  "

# Row 29 - Fetch.AI
Set-Row 29 $null $null $null "  +1.26%  "

# Row 30 - PancakeSwap
Set-Row 30 $null $null $null "  -0.54%  "

# Row 31 - Monero
Set-Row 31 $null $null "161.32" "  +2.48%  "

# Row 32 - WhiteBITCoin
Set-Row 32 $null $null "18.94" "  -0.36%  "

# Row 33 - USDe
Set-Row 33 $null $null $null "  -0.09%  "

# Row 34 - EthereumClassic
Set-Row 34 $null $null "17.77" "  +0.19%  "

# Row 35 - Kaspa
Set-Row 35 $null $null $null "  -5.94%  "

# Row 37 - RenderToken
Set-Row 37 $null $null "4.26" "  -3.30%  "

# Row 38 - Stacks
Set-Row 38 $null $null $null "  +1.16%  "

# Row 39 - ImmutableX
Set-Row 39 $null $null $null "  -1.59%  "

# Row 40 - dogwifhat
Set-Row 40 $null $null "2.01" "  -0.49%  "

# Row 41 - Filecoin
Set-Row 41 $null $null "3.34" "  +0.53%  "

# Row 42 - Aave
Set-Row 42 $null $null "128.82" "  -1.81%  "

# Row 43 - Cronos
Set-Row 43 $null $null $null "  -0.05%  "

# Row 44 - ARBITRUM
Set-Row 44 $null $null "0.477" "  -0.06%  "

# Row 45 - Mantle
Set-Row 45 $null $null $null "  -0.43%  "

# Row 46 - Stellar
Set-Row 46 $null $null "0.0914" "  +1.03%  "

# Row 47 - BitgetToken
Set-Row 47 $null $null $null "  +0.80%  "

# Row 48 - Optimism
Set-Row 48 $null $null $null "  -5.27%  "

# Row 49 - BabyDogeCoin
Set-Row 49 $null $null "0.0₆0207" "  +5.65%  "

# Row 50 - InjectiveProtocol
Set-Row 50 $null $null "16.58" "  -0.90%  "

# Row 51 - THORChain
Set-Row 51 $null $null "4.76" "  -6.22%  "
